# Auto-generated Excel COM-interop script applying the market-data refresh diff
# to the Brynhildr_Profits workbook (per-sheet leve-profit recalculation columns H-N).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 1544.5294
$ws.Range("I6").Value = 1404.2142
$ws.Range("K6").Value = 4212.642599999999
$ws.Range("M6").Value = -4100.642599999999
$ws.Range("H18").Value = 873.6667
$ws.Range("I18").Value = 873.6667
$ws.Range("K18").Value = 873.6667
$ws.Range("M18").Value = -589.6667
$ws.Range("H97").Value = 9110.385
$ws.Range("J97").Value = 12929.375
$ws.Range("L97").Value = 38788.125
$ws.Range("N97").Value = -39780.125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H14").Value = 2362.5
$ws.Range("J14").Value = 4500
$ws.Range("L14").Value = 4500
$ws.Range("N14").Value = -4850
$ws.Range("H16").Value = 148939.86
$ws.Range("I16").Value = 173513.17
$ws.Range("J16").Value = 1500
$ws.Range("K16").Value = 173513.17
$ws.Range("L16").Value = 1500
$ws.Range("M16").Value = -173226.17
$ws.Range("N16").Value = -2074
$ws.Range("H19").Value = 1335.6666
$ws.Range("I19").Value = 1503.5
$ws.Range("J19").Value = 1000
$ws.Range("K19").Value = 1503.5
$ws.Range("L19").Value = 1000
$ws.Range("M19").Value = -1274.5
$ws.Range("N19").Value = -1458
$ws.Range("H32").Value = 977628.9399999999
$ws.Range("I32").Value = 1112174
$ws.Range("J32").Value = 16592.572
$ws.Range("K32").Value = 1112174
$ws.Range("L32").Value = 16592.572
$ws.Range("M32").Value = -1111887
$ws.Range("N32").Value = -17166.572
$ws.Range("H110").Value = 2378
$ws.Range("I110").Value = 4268
$ws.Range("K110").Value = 4268
$ws.Range("M110").Value = -2223
$ws.Range("H132").Value = 5790.3076
$ws.Range("I132").Value = 3612.8333
$ws.Range("K132").Value = 10838.4999
$ws.Range("M132").Value = -8308.499899999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 4945.2144
$ws.Range("J105").Value = 8749.25
$ws.Range("L105").Value = 8749.25
$ws.Range("N105").Value = -12243.25
$ws.Range("H132").Value = 99000
$ws.Range("J132").Value = 99000
$ws.Range("L132").Value = 99000
$ws.Range("N132").Value = -109120
$ws.Range("H134").Value = 3087889.8
$ws.Range("I134").Value = 1508.75
$ws.Range("K134").Value = 4526.25
$ws.Range("M134").Value = -1991.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 749.6923
$ws.Range("J19").Value = 202
$ws.Range("L19").Value = 202
$ws.Range("N19").Value = -542
$ws.Range("H21").Value = 732.6667
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 732.6667
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 732.6667
$ws.Range("M21").ClearContents()
$ws.Range("N21").Value = -1202.6667
$ws.Range("H24").Value = 749.6923
$ws.Range("J24").Value = 202
$ws.Range("L24").Value = 202
$ws.Range("N24").Value = -542
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("M25").ClearContents()
$ws.Range("H99").Value = 20574.363
$ws.Range("I99").Value = 35453.332
$ws.Range("J99").Value = 2719.6
$ws.Range("K99").Value = 35453.332
$ws.Range("L99").Value = 2719.6
$ws.Range("M99").Value = -33955.332
$ws.Range("N99").Value = -5715.6
$ws.Range("H105").Value = 30755.715
$ws.Range("I105").Value = 20908
$ws.Range("K105").Value = 20908
$ws.Range("M105").Value = -19161
$ws.Range("H122").Value = 16390.889
$ws.Range("J122").Value = 27593.6
$ws.Range("L122").Value = 82780.79999999999
$ws.Range("N122").Value = -87680.79999999999
$ws.Range("H126").Value = 20574.363
$ws.Range("I126").Value = 35453.332
$ws.Range("J126").Value = 2719.6
$ws.Range("K126").Value = 106359.996
$ws.Range("L126").Value = 8158.799999999999
$ws.Range("M126").Value = -103889.996
$ws.Range("N126").Value = -13098.8
$ws.Range("H134").Value = 3454.111
$ws.Range("I134").Value = 3356.1924
$ws.Range("K134").Value = 10068.5772
$ws.Range("M134").Value = -7533.5772

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H99").Value = 12900
$ws.Range("I99").Value = 10375
$ws.Range("K99").Value = 31125
$ws.Range("M99").Value = -28879
$ws.Range("I107").Value = 1317.8
$ws.Range("J107").Value = 4934.4
$ws.Range("K107").Value = 3953.4
$ws.Range("L107").Value = 14803.2
$ws.Range("M107").Value = -2033.4
$ws.Range("N107").Value = -18643.2
$ws.Range("H121").Value = 14958.096
$ws.Range("J121").Value = 19399.625
$ws.Range("L121").Value = 58198.875
$ws.Range("N121").Value = -60818.875
$ws.Range("H131").Value = 6588.276
$ws.Range("I131").Value = 1884.125
$ws.Range("J131").Value = 8380.333000000001
$ws.Range("K131").Value = 5652.375
$ws.Range("L131").Value = 25140.999
$ws.Range("M131").Value = -612.375
$ws.Range("N131").Value = -35220.999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 5018750
$ws.Range("I18").Value = 10012500
$ws.Range("K18").Value = 10012500
$ws.Range("M18").Value = -10012207
$ws.Range("H21").Value = 9757.714
$ws.Range("I21").Value = 3326
$ws.Range("K21").Value = 3326
$ws.Range("M21").Value = -3153
$ws.Range("H30").Value = 9757.714
$ws.Range("I30").Value = 3326
$ws.Range("K30").Value = 3326
$ws.Range("M30").Value = -3221
$ws.Range("H80").Value = 1494.8
$ws.Range("I80").Value = 1487.5
$ws.Range("J80").Value = 1499.6666
$ws.Range("K80").Value = 1487.5
$ws.Range("L80").Value = 1499.6666
$ws.Range("M80").Value = -489.5
$ws.Range("N80").Value = -3495.6666
$ws.Range("H83").Value = 1494.8
$ws.Range("I83").Value = 1487.5
$ws.Range("J83").Value = 1499.6666
$ws.Range("K83").Value = 7437.5
$ws.Range("L83").Value = 7498.333000000001
$ws.Range("M83").Value = -2445.5
$ws.Range("N83").Value = -17482.333
$ws.Range("H102").Value = 2062.4614
$ws.Range("I102").Value = 986
$ws.Range("K102").Value = 986
$ws.Range("M102").Value = 636
$ws.Range("H122").Value = 35528.566
$ws.Range("I122").Value = 45254.305
$ws.Range("J122").Value = 3572.5715
$ws.Range("K122").Value = 135762.915
$ws.Range("L122").Value = 10717.7145
$ws.Range("M122").Value = -133312.915
$ws.Range("N122").Value = -15617.7145
$ws.Range("H132").Value = 8020.4375
$ws.Range("I132").Value = 5226.4287
$ws.Range("K132").Value = 15679.2861
$ws.Range("M132").Value = -13149.2861

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H23").Value = 20000
$ws.Range("I23").Value = 20000
$ws.Range("K23").Value = 20000
$ws.Range("M23").Value = -19770
$ws.Range("H45").Value = 45000
$ws.Range("J45").Value = 45000
$ws.Range("L45").Value = 45000
$ws.Range("N45").Value = -45814
$ws.Range("H55").Value = 1623.3793
$ws.Range("J55").Value = 1542.409
$ws.Range("L55").Value = 1542.409
$ws.Range("N55").Value = -1888.409
$ws.Range("H114").Value = 31500
$ws.Range("J114").Value = 23000
$ws.Range("L114").Value = 23000
$ws.Range("N114").Value = -31678
$ws.Range("H122").Value = 3003.3
$ws.Range("I122").Value = 2670.3333
$ws.Range("J122").Value = 6000
$ws.Range("K122").Value = 8010.999899999999
$ws.Range("L122").Value = 18000
$ws.Range("M122").Value = -5560.999899999999
$ws.Range("N122").Value = -22900
$ws.Range("H132").Value = 1237950.2
$ws.Range("I132").Value = 2567263.2
$ws.Range("J132").Value = 3588.1428
$ws.Range("K132").Value = 7701789.600000001
$ws.Range("L132").Value = 10764.4284
$ws.Range("M132").Value = -7699259.600000001
$ws.Range("N132").Value = -15824.4284

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H23").Value = 1379.5454
$ws.Range("I23").Value = 187.5
$ws.Range("K23").Value = 187.5
$ws.Range("M23").Value = 41.5
$ws.Range("H80").Value = 230300
$ws.Range("J80").Value = 230300
$ws.Range("L80").Value = 230300
$ws.Range("N80").Value = -232296
$ws.Range("H83").Value = 230300
$ws.Range("J83").Value = 230300
$ws.Range("L83").Value = 690900
$ws.Range("N83").Value = -700884
$ws.Range("H100").Value = 948.9167
$ws.Range("I100").Value = 648.7
$ws.Range("K100").Value = 1297.4
$ws.Range("M100").Value = -756.4000000000001
$ws.Range("H122").Value = 41626.105
$ws.Range("I122").Value = 1175.25
$ws.Range("K122").Value = 3525.75
$ws.Range("M122").Value = -1075.75
$ws.Range("H126").Value = 1392.1428
$ws.Range("I126").Value = 1133.3334
$ws.Range("K126").Value = 3400.0002
$ws.Range("M126").Value = -930.0001999999999
$ws.Range("H132").Value = 5748960
$ws.Range("I132").Value = 6946065
$ws.Range("K132").Value = 20838195
$ws.Range("M132").Value = -20835665
$ws.Range("H136").Value = 11776973
$ws.Range("I136").Value = 7248126.5
$ws.Range("K136").Value = 21744379.5
$ws.Range("M136").Value = -21741829.5
